# Reposition the five "5-Point Star" shapes on slide 3 (moving them down and
# to the left, per the "updated new RR again" commit).
#
# The target coordinates are expressed in EMU (English Metric Units), as
# stored in the underlying OOXML <a:off>, while PowerPoint's COM object
# model (Shape.Left / Shape.Top) works in points and stores them as
# single-precision floats. Converting EMU -> points -> (float32) -> EMU can
# truncate the last EMU off due to float32 rounding, so a tiny epsilon
# (half an EMU, expressed in points) is added before assigning to land back
# on the exact target EMU value after the round-trip.

$EMU_PER_POINT = 914400 / 72
$EPS_POINTS = 0.5 / $EMU_PER_POINT

function EmuToPoints($emu) {
    return ($emu / $EMU_PER_POINT) + $EPS_POINTS
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Map of shape name -> new [x,y] offset in EMU, matching the target OOXML.
$moves = @{
    "5-Point Star 6"  = @(3075925, 6987693)
    "5-Point Star 23" = @(3802197, 6982126)
    "5-Point Star 24" = @(4544848, 6969903)
    "5-Point Star 25" = @(5271181, 6982126)
    "5-Point Star 26" = @(5996543, 6982126)
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($moves.ContainsKey($shape.Name)) {
        $target = $moves[$shape.Name]
        $shape.Left = EmuToPoints $target[0]
        $shape.Top = EmuToPoints $target[1]
    }
}
